$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Final (target) values for columns D, J, K, L, M, P across rows 2-7 and 10.
# These were derived by re-ordering the original rows' data (a cyclic shuffle
# of the weekly price records) as described by the commit diff.
$rows = @{
    2  = @{ D = 44893; J = 3300; K = 1200; L = 1300; M = 1261; P = 1261 }
    3  = @{ D = 44537; J = 800;  K = 1300; L = 1400; M = 1350; P = 1350 }
    4  = @{ D = 44895; J = 200;  K = 1200; L = 1300; M = 1255; P = 1255 }
    5  = @{ D = 44883; J = 290;  K = 1400; L = 1500; M = 1434; P = 1434 }
    6  = @{ D = 44200; J = 1500; K = 1400; L = 1500; M = 1450; P = 1450 }
    7  = @{ D = 44638; J = 800;  K = 2500; L = 2800; M = 2650; P = 2650 }
    10 = @{ D = 44210; J = 1450; K = 1600; L = 1700; M = 1650; P = 1650 }
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    $ws.Range("D$r").Value = $vals.D
    $ws.Range("J$r").Value = $vals.J
    $ws.Range("K$r").Value = $vals.K
    $ws.Range("L$r").Value = $vals.L
    $ws.Range("M$r").Value = $vals.M
    $ws.Range("P$r").Value = $vals.P
}
